# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1!A1 : update the two "Binance" conversion bullet lines ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $hoja1.Range("A1")
$text = $cellA1.Value2

$lines = $text -split "`n"
for ($i = 0; $i -lt $lines.Length; $i++) {
    if ($lines[$i].Contains("1000 Bs = 1.7 = 6212.77 pesos")) {
        $lines[$i] = $lines[$i].Replace("1000 Bs = 1.7 = 6212.77 pesos", "1000 Bs = 1.69 = 6148.65 pesos")
    }
    if ($lines[$i].Contains("6212.77 pesos = 1.69 = 931.07 Bs")) {
        $lines[$i] = $lines[$i].Replace("6212.77 pesos = 1.69 = 931.07 Bs", "6148.65 pesos = 1.68 = 933.86 Bs")
    }
}
$cellA1.Value = ($lines -join "`n")

# --- tasas!N10/O10/N12/O12 : updated rate figures ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 592
$tasas.Range("O10").Value = 3640
$tasas.Range("N12").Value = 3660.78
$tasas.Range("O12").Value = 556
